$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("S3").Value = 3.8
$ws.Range("T3").Value = 1.25

# Row 4
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3.2

# Row 5
$ws.Range("G5").Value = 5
$ws.Range("I5").Value = 1.6
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 12
$ws.Range("Y5").Value = 1.8
$ws.Range("Z5").Value = 1.91
$ws.Range("AB5").Value = 26
$ws.Range("AG5").Value = 12
$ws.Range("AL5").Value = 7.5
$ws.Range("AM5").Value = 8
$ws.Range("AQ5").Value = 23

# Row 6
$ws.Range("G6").Value = 2.45
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 2.88
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.5
$ws.Range("AB6").Value = 12
$ws.Range("AC6").Value = 9.5
$ws.Range("AF6").Value = 26
$ws.Range("AN6").Value = 11
